# Update "Example1" schedule sheet with the latest codes pulled from Google Drive
# (source data snapshot: 20210914)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 4-7 (Tue/Wed/Thu/Fri schedule entries) - only one
# teacher/day/class/subject/level row remains below the header row.
$ws.Range("A4:E7").ClearContents()

# Row 1: teacher name / role label / year
$ws.Range("A1").Value = "ทดสอบ"
$ws.Range("B1").Value = "คุณครู"
$ws.Range("C1").Value = 2561

# Row 2: column headers
$ws.Range("A2").Value = "date"
$ws.Range("B2").Value = "เวลาในการสอน"
$ws.Range("C2").Value = "ชั้น"
$ws.Range("D2").Value = "รหัสวิชา"
$ws.Range("E2").Value = "ระดับชั้นเรียน"

# Row 3: the single remaining schedule entry
$ws.Range("A3").Value = "วันพุธ"
$ws.Range("B3").Value = "14:00 - 15:00"
$ws.Range("C3").Value = "ป.1/1"
$ws.Range("D3").Value = "ค 11101"
$ws.Range("E3").Value = "ประถมศึกษา"

# Match the selection left behind in the authored file
$ws.Range("K12").Select()
